# Diary.xlsx update:
#  - Rows 9-12 (w/c 18/12, 25/12, 01/01, 08/01) get their "Timeline Stage"
#    (column B) set to "Format & General research", matching the previous
#    weeks' entries.
#  - Row 10 (Christmas week) gets a short commit message.
#  - Row 12 gets a commit message about evaluating the project idea
#    shortlist.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B9").Value  = "Format & General research"
$ws.Range("B10").Value = "Format & General research"
$ws.Range("B11").Value = "Format & General research"
$ws.Range("B12").Value = "Format & General research"

$ws.Range("C10").Value = "Literally Christmas. Did nothing."
$ws.Range("C12").Value = "Evaluation of project idea shortlist. Further research, and evaluation of research."

# Reflect the scrolled-down view / new selection left by the edit.
$ws.Activate()
$ws.Range("C11").Select()
